# Update the quarterly income statement:
#  - Drop the oldest quarter (column D, "Q2 ending 1399/06")
#  - Shift all remaining quarters one column to the left
#  - Append the newest quarter (column M, "Q4 ending 1401/12")
#  - Re-state column I ("Q3 ending 1401/09") with figures produced by the
#    updated read_price algorithm (its publish date also moved to 1402-02-30 (8))

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the oldest quarter column; everything to the right shifts left.
$ws.Range("D:D").Delete() | Out-Null

# The shift leaves column M blank; clone column L's formatting into it
# (column widths follow a Q4-is-wider pattern, handled separately below).
$ws.Range("L1:L28").Copy() | Out-Null
$ws.Range("M1:M28").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# Column E is now a "quarter ending in Q4" column (width 31); mirror that
# width onto the newly appended column M, which is also a Q4 column now.
$ws.Columns("M").ColumnWidth = $ws.Columns("E").ColumnWidth

# New period header + publish date for the newly appended quarter.
$ws.Range("M8").Value = "فصل چهارم منتهی به 1401/12"
$ws.Range("M9").Value = "1402-02-30"

# The republished Q3 1401/09 figures now show this later publish date.
$ws.Range("I9").Value = "1402-02-30 (8)"

# New quarter (column M) financial figures.
$ws.Range("M11").Value = 40278228
$ws.Range("M12").Value = -27963235
$ws.Range("M13").Value = 12314993
$ws.Range("M14").Value = -11337419
$ws.Range("M15").Value = 0
$ws.Range("M16").Value = 5840423
$ws.Range("M17").Value = 6817997
$ws.Range("M18").Value = 0
$ws.Range("M19").Value = -2114913
$ws.Range("M20").Value = 4703084
$ws.Range("M21").Value = 5383870
$ws.Range("M22").Value = 10086954
$ws.Range("M23").Value = 0
$ws.Range("M24").Value = 10086954
$ws.Range("M25").Value = 4203
$ws.Range("M26").Value = 2400000
$ws.Range("M27").Value = 4203

# Column I (Q3 1401/09) figures recomputed by the new read_price algorithm.
$ws.Range("I12").Value = -32855503
$ws.Range("I13").Value = 5378444
$ws.Range("I14").Value = -8541335
$ws.Range("I17").Value = -4004444
$ws.Range("I20").Value = -2864948
$ws.Range("I22").Value = -2864948
$ws.Range("I24").Value = -2864948
$ws.Range("I25").Value = -1194
$ws.Range("I27").Value = -1194
